$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" summary text (cell A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.1 = 7939.24 pesos`n✅ 7939.24 pesos = 2.09 = 954.97 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update the N10/O10/N12/O12 rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 476.5
$ws2.Range("O10").Value = 3783.05
$ws2.Range("N12").Value = 3791
$ws2.Range("O12").Value = 456
